$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '28.044.30', '  -0.18%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.870.00', '  -1.22%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.003', '  +0.26%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '312.65', '  -0.51%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.003', '  +0.24%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5081', '  +1.10%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3809', '  -2.24%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08312', '  -9.85%  '),
    @(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.111', '  -1.63%  '),
    @(11, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.214', '  -2.61%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.870.06', '  -1.25%  '),
    @(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.48', '  -1.64%  '),
    @(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.188', '  -1.51%  '),
    @(15, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.003', '  +0.21%  '),
    @(16, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001096', '  -1.05%  '),
    @(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '90.72', '  -1.41%  '),
    @(18, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06630', '  -0.10%  '),
    @(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.86', '  +0.08%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  +0.19%  '),
    @(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.012', '  -3.46%  '),
    @(22, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '28.080.82', '  -0.26%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.11', '  -2.35%  '),
    @(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.263', '  -2.59%  '),
    @(25, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.578', '  +1.46%  '),
    @(26, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.086.66', '  -1.22%  '),
    @(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '157.26', '  -0.80%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '20.53', '  -1.41%  '),
    @(29, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '125.97', '  -0.81%  '),
    @(30, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1057', '  +0.24%  '),
    @(31, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.043', '  -3.06%  '),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.596', '  -0.10%  '),
    @(33, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.590', '  -0.50%  '),
    @(34, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.735', '  +2.65%  '),
    @(35, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02443', '  +1.45%  '),
    @(36, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06562', '  -0.68%  '),
    @(37, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2164', '  -1.57%  '),
    @(38, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.209', '  -0.88%  '),
    @(39, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6446', '  -0.11%  '),
    @(40, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.237', '  -7.89%  '),
    @(41, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.30', '  -2.37%  '),
    @(42, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.879', '  -1.45%  '),
    @(43, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6139', '  +1.29%  '),
    @(44, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.97', '  -2.73%  '),
    @(45, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.289', '  -0.76%  '),
    @(46, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.663', '  -0.85%  '),
    @(47, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.004', '  -0.10%  '),
    @(48, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.213', '  +1.63%  '),
    @(49, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '121.28', '  -0.41%  '),
    @(50, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '80.15', '  +1.34%  '),
    @(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06864', '  -0.49%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $priceCell = $ws.Cells.Item($r, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
}
